$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.037.10"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "3.009.91"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.71"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +5.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.73"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +4.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +3.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.56"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +6.09%  "
$ws.Range("E10").Value = "  +6.39%  "
$ws.Range("E11").Value = "  +3.19%  "
$ws.Range("E12").Value = "  +2.60%  "
$ws.Range("D13").Value = "3.526.72"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.76"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("E15").Value = "  +11.77%  "
$ws.Range("D16").Value = "57.059.38"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "3.010.86"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.97"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +5.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.62"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.91"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.14"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +2.93%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +4.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.70"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +5.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +6.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "0.0₃0918"
$ws.Range("E27").Value = "  +8.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.68"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +2.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.15"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +7.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +6.27%  "
$ws.Range("E31").Value = "  +6.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.70"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +5.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.81"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +4.93%  "
$ws.Range("E34").Value = "  +4.51%  "
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("E37").Value = "  +3.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.25"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +3.13%  "
$ws.Range("D39").Value = "3.043.84"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.25"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").Value = "2.294.27"
$ws.Range("E42").Value = "  +7.98%  "
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.72"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +4.81%  "
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.99"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +8.52%  "
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.88"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +5.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.33"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0878"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +4.76%  "
